# Apply "Added team record to data" edit:
# Add three new columns (Wins, Losses, Ties) to the right of the existing
# data, with the same header formatting as the other header cells, and
# populate every data row (2-43) with the team record values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from an existing header cell (A1) onto the
# three new header cells so they match the rest of row 1 (bold, centered,
# bordered).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set the new header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate the team record for each data row (rows 2 through 43).
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 52
    $ws.Cells.Item($r, 31).Value = 62
    $ws.Cells.Item($r, 32).Value = 0
}
